$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new property row (38) describing the new `checkDuplicateTitle`
# boolean property, mirroring the layout of the existing property rows
# (label in column B, comment/description in column D).
$ws.Range("B38").Value = " boolean checkDuplicateTitle"
$ws.Range("D38").Value = " //타이틀 중복체크 수행여부."

# Keep the current selection on the sheet, same as the authored workbook.
$ws.Range("B24").Select() | Out-Null
